# Update the arithmetic drill table: replace each "two-digit divided by
# one-digit" problem with the new problem from the commit.
#
# Replacements are applied by (row, column) position in the single table
# rather than by text search, because a couple of the old values
# ("36÷6=") are duplicated in the sheet but map to *different* new
# values depending on which cell they're in.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (problems 1-5)
$t.Cell(1,1).Range.Text = "35÷3="
$t.Cell(1,2).Range.Text = "32÷8="
$t.Cell(1,3).Range.Text = "45÷4="
$t.Cell(1,4).Range.Text = "55÷6="
$t.Cell(1,5).Range.Text = "59÷8="

# Row 5 (problems 6-10)
$t.Cell(5,1).Range.Text = "92÷3="
$t.Cell(5,2).Range.Text = "78÷3="
$t.Cell(5,3).Range.Text = "97÷2="
$t.Cell(5,4).Range.Text = "19÷7="
$t.Cell(5,5).Range.Text = "46÷2="

# Row 9 (problems 11-15)
$t.Cell(9,1).Range.Text = "70÷2="
$t.Cell(9,2).Range.Text = "32÷3="
$t.Cell(9,3).Range.Text = "77÷4="
$t.Cell(9,4).Range.Text = "53÷8="
$t.Cell(9,5).Range.Text = "96÷7="

# Row 13 (problems 16-20)
$t.Cell(13,1).Range.Text = "35÷2="
$t.Cell(13,2).Range.Text = "41÷9="
$t.Cell(13,3).Range.Text = "88÷3="
$t.Cell(13,4).Range.Text = "74÷7="
$t.Cell(13,5).Range.Text = "91÷2="

# Row 17 (problems 21-25)
$t.Cell(17,1).Range.Text = "20÷8="
$t.Cell(17,2).Range.Text = "17÷5="
$t.Cell(17,3).Range.Text = "33÷9="
$t.Cell(17,4).Range.Text = "32÷8="
$t.Cell(17,5).Range.Text = "45÷4="
